$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow")

# Update rows 2-9 of the Flow sheet with new values (ActivityId, ActivityIdTarget, Probability, DurationIdle)
$ws.Range("A2").Value = "<Start>"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 5

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 15

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "3a"
$ws.Range("C4").Value = 0.3
$ws.Range("D4").Value = 0

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "3b"
$ws.Range("C5").Value = 0.6
$ws.Range("D5").Value = 10

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "<End>"
$ws.Range("C6").Value = 0.1
$ws.Range("D6").Value = 0

$ws.Range("A7").Value = "3a"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0

$ws.Range("A8").Value = "3b"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 10

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "<End>"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0

# Remove row 10 entirely (shift cells up)
$ws.Range("A10:D10").Delete()

# Remove conditional formatting rules
$ws.Cells.FormatConditions.Delete()

# Update selection to A9
$ws.Range("A9").Select()
